# Auto-generated Excel COM-interop script applying the Odin_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I132").Value = 313309.12
$ws.Range("K132").Value = 939927.36
$ws.Range("M132").Value = -937397.36
$ws.Range("H132").Value = 282745.94

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J13").Value = 38000
$ws.Range("N13").Value = -38336
$ws.Range("H13").Value = 38000
$ws.Range("L13").Value = 38000
$ws.Range("N50").Value = -65148
$ws.Range("H50").Value = 64000
$ws.Range("L50").Value = 64000
$ws.Range("J50").Value = 64000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N16").Value = -13240.333
$ws.Range("L16").Value = 12666.333
$ws.Range("K16").Value = 71433976
$ws.Range("I16").Value = 71433976
$ws.Range("M16").Value = -71433689
$ws.Range("J16").Value = 12666.333
$ws.Range("H16").Value = 50007584
$ws.Range("N20").Value = -140472
$ws.Range("H20").Value = 140000
$ws.Range("J20").Value = 140000
$ws.Range("L20").Value = 140000
$ws.Range("L30").Value = 140000
$ws.Range("H30").Value = 140000
$ws.Range("N30").Value = -140182
$ws.Range("J30").Value = 140000
$ws.Range("H43").Value = 29424.75
$ws.Range("N43").Value = -29792.75
$ws.Range("J43").Value = 29424.75
$ws.Range("L43").Value = 29424.75
$ws.Range("H68").Value = 35000
$ws.Range("N68").Value = -36498
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("K69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("H70").Value = 45000
$ws.Range("K70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("L71").Value = 105000
$ws.Range("H71").Value = 35000
$ws.Range("N71").Value = -112488
$ws.Range("J71").Value = 35000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 45000
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("I73").Value = 0
$ws.Range("H74").Value = 24387.4
$ws.Range("J74").Value = 24387.4
$ws.Range("N74").Value = -26135.4
$ws.Range("L74").Value = 24387.4
$ws.Range("J77").Value = 24387.4
$ws.Range("N77").Value = -81898.20000000001
$ws.Range("H77").Value = 24387.4
$ws.Range("L77").Value = 73162.20000000001
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15812
$ws.Range("J88").Value = 15000
$ws.Range("H88").Value = 15000
$ws.Range("L91").Value = 15000
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("N91").Value = -17808
$ws.Range("L92").Value = 19312
$ws.Range("H92").Value = 19312
$ws.Range("J92").Value = 19312
$ws.Range("N92").Value = -24304
$ws.Range("J93").Value = 15000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -18744
$ws.Range("M93").Value = -28128
$ws.Range("K93").Value = 30000
$ws.Range("H93").Value = 22500
$ws.Range("I93").Value = 30000
$ws.Range("L95").Value = 13899.714
$ws.Range("H95").Value = 13899.714
$ws.Range("J95").Value = 13899.714
$ws.Range("N95").Value = -19391.714
$ws.Range("L97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("I99").Value = 10103845
$ws.Range("K99").Value = 10103845
$ws.Range("H99").Value = 5558888.5
$ws.Range("M99").Value = -10102347
$ws.Range("J101").Value = 29424.75
$ws.Range("H101").Value = 29424.75
$ws.Range("L101").Value = 29424.75
$ws.Range("N101").Value = -35914.75
$ws.Range("M102").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("I104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("K106").Value = 0
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("H107").Value = 690
$ws.Range("M107").Value = 1230
$ws.Range("K107").Value = 690
$ws.Range("I107").Value = 690
$ws.Range("L108").Value = 0
$ws.Range("H108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("J108").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N113").Value = -17006.333
$ws.Range("I113").Value = 71433976
$ws.Range("M113").Value = -71431806
$ws.Range("H113").Value = 50007584
$ws.Range("K113").Value = 71433976
$ws.Range("J113").Value = 12666.333
$ws.Range("L113").Value = 12666.333
$ws.Range("N117").Value = -89173
$ws.Range("J117").Value = 79995
$ws.Range("L117").Value = 79995
$ws.Range("H117").Value = 79995
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("N120").ClearContents()
$ws.Range("L120").Value = 0
$ws.Range("H120").Value = 0
$ws.Range("L121").Value = 116312.5
$ws.Range("H121").Value = 116312.5
$ws.Range("J121").Value = 116312.5
$ws.Range("N121").Value = -118932.5
$ws.Range("H122").Value = 12332
$ws.Range("I122").Value = 2655.5
$ws.Range("K122").Value = 7966.5
$ws.Range("J122").Value = 16202.6
$ws.Range("M122").Value = -5516.5
$ws.Range("N122").Value = -53507.8
$ws.Range("L122").Value = 48607.8
$ws.Range("N124").ClearContents()
$ws.Range("J124").Value = 0
$ws.Range("H124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 66666.664
$ws.Range("N125").Value = -71586.664
$ws.Range("L125").Value = 66666.664
$ws.Range("J125").Value = 66666.664
$ws.Range("K126").Value = 30311535
$ws.Range("M126").Value = -30309065
$ws.Range("I126").Value = 10103845
$ws.Range("H126").Value = 5558888.5
$ws.Range("N128").Value = -149960
$ws.Range("L128").Value = 140000
$ws.Range("J128").Value = 140000
$ws.Range("H128").Value = 140000
$ws.Range("M134").Value = -166674861
$ws.Range("K134").Value = 166677396
$ws.Range("I134").Value = 55559132
$ws.Range("H134").Value = 43485164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M2").Value = -5793.25
$ws.Range("N2").Value = -940
$ws.Range("J2").Value = 119
$ws.Range("H2").Value = 465.15
$ws.Range("L2").Value = 714
$ws.Range("I2").Value = 984.375
$ws.Range("K2").Value = 5906.25
$ws.Range("K43").Value = 0
$ws.Range("H43").Value = 5000
$ws.Range("N43").Value = -15228
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5000
$ws.Range("M43").ClearContents()
$ws.Range("L43").Value = 15000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L21").Value = 21326.666
$ws.Range("H21").Value = 21326.666
$ws.Range("N21").Value = -21672.666
$ws.Range("J21").Value = 21326.666
$ws.Range("I22").Value = 5006
$ws.Range("H22").Value = 5006
$ws.Range("K22").Value = 5006
$ws.Range("M22").Value = -4477
$ws.Range("L30").Value = 21326.666
$ws.Range("H30").Value = 21326.666
$ws.Range("N30").Value = -21536.666
$ws.Range("J30").Value = 21326.666
$ws.Range("J49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("L49").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("I52").Value = 0
$ws.Range("M53").Value = -3369
$ws.Range("I53").Value = 4000
$ws.Range("K53").Value = 4000
$ws.Range("H53").Value = 4000
$ws.Range("K55").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("I55").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("N57").Value = -34973
$ws.Range("H57").Value = 33333
$ws.Range("L57").Value = 33333
$ws.Range("J57").Value = 33333
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("I58").Value = 30041
$ws.Range("N58").Value = -30554
$ws.Range("M58").Value = -29764
$ws.Range("H58").Value = 30020.5
$ws.Range("J58").Value = 30000
$ws.Range("K58").Value = 30041
$ws.Range("L58").Value = 30000
$ws.Range("H62").Value = 49999.75
$ws.Range("M62").Value = -49313
$ws.Range("I62").Value = 49999
$ws.Range("K62").Value = 49999
$ws.Range("H65").Value = 49999.75
$ws.Range("K65").Value = 149997
$ws.Range("M65").Value = -146565
$ws.Range("I65").Value = 49999
$ws.Range("I80").Value = 5288.5
$ws.Range("K80").Value = 5288.5
$ws.Range("H80").Value = 7916.8335
$ws.Range("M80").Value = -4290.5
$ws.Range("K83").Value = 26442.5
$ws.Range("I83").Value = 5288.5
$ws.Range("M83").Value = -21450.5
$ws.Range("H83").Value = 7916.8335
